# Meridian Development Group Master workbook - import field backfill
#
# Adds columns that were missing from several import handlers:
#   - RFIs:        cost_impact, schedule_impact_days
#   - Submittals:  submittal_number (new first column), review_comments (new last column)
#   - Maintenance: actual_cost, notes
#   - Equipment:   last_maintenance_date, next_maintenance_date

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    # Force a value to be stored as text even when it looks like a number
    # or a date (Excel would otherwise auto-convert it on assignment).
    if ($null -eq $text) {
        $text = ""
    }
    $ws.Cells.Item($row, $col).Value2 = "'" + $text
}

# ---------------------------------------------------------------------------
# RFIs (sheet "RFIs") - add cost_impact (G) and schedule_impact_days (H)
# ---------------------------------------------------------------------------
$rfis = $wb.Worksheets.Item("RFIs")

$rfis.Cells.Item(1, 7).Value2 = "cost_impact"
$rfis.Cells.Item(1, 8).Value2 = "schedule_impact_days"

$rfiData = @(
    @("12500", "3"),
    @("8200", "0"),
    @("35000", "7"),
    @("0", "5"),
    @("15000", "0"),
    @("0", "0"),
    @("22000", "2"),
    @("18500", "0"),
    @("6800", "5"),
    @("4200", "3"),
    @("2100", "0"),
    @("8500", "2"),
    @("12000", "5"),
    @("9800", "3"),
    @("3500", "0"),
    @("5600", "2")
)

for ($i = 0; $i -lt $rfiData.Length; $i++) {
    $row = $i + 2
    Set-TextCell $rfis $row 7 $rfiData[$i][0]
    Set-TextCell $rfis $row 8 $rfiData[$i][1]
}

Write-Host "RFIs updated:" $rfis.UsedRange.Address()

# ---------------------------------------------------------------------------
# Submittals (sheet "Submittals")
#   - insert submittal_number as the new first column (A), shifting the rest
#     of the columns (title..status) one place to the right (B..F)
#   - append review_comments as the new last column (G)
# ---------------------------------------------------------------------------
$submittals = $wb.Worksheets.Item("Submittals")

$submittals.Columns("A").Insert()

$submittals.Cells.Item(1, 1).Value2 = "submittal_number"
$submittals.Cells.Item(1, 7).Value2 = "review_comments"

$submittalData = @(
    @("SUB-001", "Approved with comments - verify connection details at moment frame locations per structural addendum 3"),
    @("SUB-002", "Approved - mockup test passed. Proceed with fabrication per approved sample"),
    @("SUB-003", "Approved - equipment selections meet spec requirements. Verify electrical connections with Division 26"),
    @("SUB-004", "Approved as noted - confirm arc flash labeling per NFPA 70E"),
    @("SUB-005", ""),
    @("SUB-006", ""),
    @("SUB-007", ""),
    @("SUB-008", ""),
    @("SUB-009", "Approved - mill certificates provided. Pile driving contractor to submit driving criteria separately"),
    @("SUB-010", "Approved with conditions - maintain 28-day cylinder break results above 6000 PSI"),
    @("SUB-011", ""),
    @("SUB-012", ""),
    @("SUB-013", "Approved - bearing capacity and movement range verified by structural engineer"),
    @("SUB-014", "Approved as submitted"),
    @("SUB-015", "Approved - HPC mix meets AASHTO requirements for bridge deck application"),
    @("SUB-016", "Approved with note - torque values per manufacturer ICC-ES report"),
    @("SUB-017", "Approved - RF shielding and vibration isolation requirements confirmed with vendor"),
    @("SUB-018", "Approved - Type K copper per RFI clarification. ASSE 6010 certification required for installer"),
    @("SUB-019", "Approved - NEC 700 and NFPA 110 compliance verified"),
    @("SUB-020", "Approved - Florida Building Code NOA documentation provided"),
    @("SUB-021", "Approved with comments - verify panel connection hardware per structural detail SD-12"),
    @("SUB-022", "Approved - health department flow rate requirements met")
)

for ($i = 0; $i -lt $submittalData.Length; $i++) {
    $row = $i + 2
    $submittals.Cells.Item($row, 1).Value2 = $submittalData[$i][0]
    $comment = $submittalData[$i][1]
    if ($comment -eq "") {
        $submittals.Cells.Item($row, 7).NumberFormat = "@"
        $submittals.Cells.Item($row, 7).Formula = ""
    } else {
        $submittals.Cells.Item($row, 7).Value2 = $comment
    }
}

Write-Host "Submittals updated:" $submittals.UsedRange.Address()

# ---------------------------------------------------------------------------
# Maintenance (sheet "Maintenance") - append actual_cost (I) and notes (J)
# ---------------------------------------------------------------------------
$maintenance = $wb.Worksheets.Item("Maintenance")

$maintenance.Cells.Item(1, 9).Value2 = "actual_cost"
$maintenance.Cells.Item(1, 10).Value2 = "notes"

$maintenanceData = @(
    @("4200", "Compressor replaced same day. Tenant satisfied with response time."),
    @("", "Weekend work required to minimize tenant impact."),
    @("", "Coordinate with elevator vendor for after-hours access."),
    @("3100", "Upgraded to heat pump model. Slight cost increase but better long-term efficiency."),
    @("1850", "Motor replaced and pool reopened within 48 hours."),
    @("", "Roofer identified failed flashing at parapet wall. Temporary repair in place."),
    @("850", "Grease trap cleaned and certified. Next service due October 2025."),
    @("", "Glass on order. 3-week lead time from manufacturer."),
    @("1450", "Repaired and operational. Recommended preventive service for all 6 bay doors."),
    @("", "Schedule with fire protection contractor. Must complete before insurance renewal.")
)

for ($i = 0; $i -lt $maintenanceData.Length; $i++) {
    $row = $i + 2
    $cost = $maintenanceData[$i][0]
    if ($cost -eq "") {
        $maintenance.Cells.Item($row, 9).NumberFormat = "@"
        $maintenance.Cells.Item($row, 9).Formula = ""
    } else {
        Set-TextCell $maintenance $row 9 $cost
    }
    $maintenance.Cells.Item($row, 10).Value2 = $maintenanceData[$i][1]
}

Write-Host "Maintenance updated:" $maintenance.UsedRange.Address()

# ---------------------------------------------------------------------------
# Equipment (sheet "Equipment") - append last_maintenance_date (J) and
# next_maintenance_date (K)
# ---------------------------------------------------------------------------
$equipment = $wb.Worksheets.Item("Equipment")

$equipment.Cells.Item(1, 10).Value2 = "last_maintenance_date"
$equipment.Cells.Item(1, 11).Value2 = "next_maintenance_date"

$equipmentData = @(
    @("2025-02-10", "2025-08-10"),
    @("2025-01-20", "2026-01-20"),
    @("2025-04-10", "2026-04-10"),
    @("2025-06-22", "2025-12-22"),
    @("2025-09-01", "2026-01-01"),
    @("2025-03-01", "2025-09-01"),
    @("2025-03-15", "2025-09-15"),
    @("2025-05-01", "2025-08-01"),
    @("2025-07-10", "2026-07-10"),
    @("2025-04-28", "2025-10-28"),
    @("2025-06-15", "2026-06-15"),
    @("2025-08-20", "2026-08-20")
)

for ($i = 0; $i -lt $equipmentData.Length; $i++) {
    $row = $i + 2
    Set-TextCell $equipment $row 10 $equipmentData[$i][0]
    Set-TextCell $equipment $row 11 $equipmentData[$i][1]
}

Write-Host "Equipment updated:" $equipment.UsedRange.Address()
